{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Remove the first two paragraphs: the \"Wednesday 10/1/2025 12 PM ET\"\n// date/time header and the \"Hurricane Imelda\" title paragraph.\nparagraphs.items[0].delete();\nparagraphs.items[1].delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the first two paragraphs: the \"Wednesday 10/1/2025 12 PM ET\"\n# date/time header and the \"Hurricane Imelda\" title paragraph.\n$d.Paragraphs(1).Range.Delete()\n$d.Paragraphs(1).Range.Delete()\n"}
